# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Octubre de 2020 a las 06:43"

# Row 5 - India
$ws.Range("B5").Value = 6906151
$ws.Range("C5").Value = 2339
$ws.Range("D5").Value = 5906069
$ws.Range("E5").Value = 893561

# Row 58 - Uzbekistan
$ws.Range("B58").Value = 60122
$ws.Range("C58").Value = 96
$ws.Range("D58").Value = 56976
$ws.Range("E58").Value = 2650

# Row 142 - Tailandia
$ws.Range("B142").Value = 3628
$ws.Range("C142").Value = 6
$ws.Range("D142").Value = 3441
$ws.Range("E142").Value = 128

# Rows 154/155 - Letonia and Belice swap places (row 154 becomes Belice,
# row 155 becomes Letonia), each carrying updated figures.
$ws.Range("A154").Value = "Belice"
$ws.Range("B154").Value = 2373
$ws.Range("C154").Value = 63
$ws.Range("D154").Value = 1459
$ws.Range("E154").Value = 880
$ws.Range("H154").Value = 34

$ws.Range("A155").Value = "Letonia"
$ws.Range("B155").Value = 2370
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 1322
$ws.Range("E155").Value = 1008
$ws.Range("H155").Value = 40

# Row 187 - Butan
$ws.Range("D187").Value = 255
$ws.Range("E187").Value = 49
